$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '68.437.64'
$c.Style = $s
$ws.Range('E2').Value = '  +1.20%  '
$c = $ws.Range('D3')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '3.748.38'
$c.Style = $s
$ws.Range('E3').Value = '  -0.70%  '
$ws.Range('E4').Value = '  +0.02%  '
$c = $ws.Range('D5')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '595.27'
$c.Style = $s
$ws.Range('E5').Value = '  -0.17%  '
$c = $ws.Range('D6')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '166.96'
$c.Style = $s
$ws.Range('E6').Value = '  -0.84%  '
$c = $ws.Range('D7')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '3.749.01'
$c.Style = $s
$ws.Range('E7').Value = '  -0.64%  '
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('E9').Value = '  -0.74%  '
$ws.Range('E10').Value = '  -2.57%  '
$c = $ws.Range('D11')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '6.46'
$c.Style = $s
$ws.Range('E11').Value = '  -0.88%  '
$c = $ws.Range('D12')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.446'
$c.Style = $s
$ws.Range('E12').Value = '  -1.10%  '
$ws.Range('E13').Value = '  -5.92%  '
$c = $ws.Range('D14')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '35.98'
$c.Style = $s
$ws.Range('E14').Value = '  -0.79%  '
$c = $ws.Range('D15')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '4.378.96'
$c.Style = $s
$ws.Range('E15').Value = '  -0.77%  '
$c = $ws.Range('D16')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '3.756.28'
$c.Style = $s
$ws.Range('E16').Value = '  -0.66%  '
$c = $ws.Range('D17')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '68.461.29'
$c.Style = $s
$ws.Range('E17').Value = '  +1.26%  '
$ws.Range('E18').Value = '  -2.91%  '
$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Range('D19')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.112'
$c.Style = $s
$ws.Range('E19').Value = '  -0.08%  '
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range('D20')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '7.00'
$c.Style = $s
$ws.Range('E20').Value = '  -2.48%  '
$c = $ws.Range('D21')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '10.75'
$c.Style = $s
$ws.Range('E21').Value = '  +2.58%  '
$c = $ws.Range('D22')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '465.26'
$c.Style = $s
$ws.Range('E22').Value = '  -0.29%  '
$c = $ws.Range('D23')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.697'
$c.Style = $s
$ws.Range('E23').Value = '  -2.59%  '
$c = $ws.Range('D24')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '84.52'
$c.Style = $s
$ws.Range('E24').Value = '  +1.00%  '
$ws.Range('E25').Value = '  -1.89%  '
$ws.Range('E26').Value = '  -0.23%  '
$c = $ws.Range('D27')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '11.98'
$c.Style = $s
$ws.Range('E27').Value = '  -0.98%  '
$ws.Range('E28').Value = '  -0.05%  '
$c = $ws.Range('D29')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '10.01'
$c.Style = $s
$ws.Range('E29').Value = '  -2.92%  '
$c = $ws.Range('D30')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '3.895.53'
$c.Style = $s
$ws.Range('E30').Value = '  -0.88%  '
$ws.Range('E31').Value = '  -4.33%  '
$c = $ws.Range('D32')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '7.28'
$c.Style = $s
$ws.Range('E32').Value = '  -4.12%  '
$c = $ws.Range('D33')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '29.79'
$c.Style = $s
$ws.Range('E33').Value = '  -2.28%  '
$ws.Range('E34').Value = '  -2.00%  '
$c = $ws.Range('D35')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '9.20'
$c.Style = $s
$ws.Range('E35').Value = '  +1.03%  '
$c = $ws.Range('D37')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '3.704.03'
$c.Style = $s
$ws.Range('E37').Value = '  -1.07%  '
$ws.Range('E38').Value = '  -2.55%  '
$ws.Range('E39').Value = '  -8.90%  '
$ws.Range('E40').Value = '  +0.93%  '
$c = $ws.Range('D41')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = $s
$ws.Range('E41').Value = '  -0.29%  '
$c = $ws.Range('D42')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '5.81'
$c.Style = $s
$ws.Range('E42').Value = '  +0.30%  '
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('E45').Value = '  -2.09%  '
$c = $ws.Range('D46')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '43.79'
$c.Style = $s
$ws.Range('E46').Value = '  +11.81%  '
$c = $ws.Range('D47')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '8.55'
$c.Style = $s
$ws.Range('E47').Value = '  -0.97%  '
$c = $ws.Range('D49')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '45.93'
$c.Style = $s
$ws.Range('E49').Value = '  +0.48%  '
$c = $ws.Range('D50')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '146.29'
$c.Style = $s
$ws.Range('E50').Value = '  +4.08%  '
$c = $ws.Range('D51')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '389.32'
$c.Style = $s
$ws.Range('E51').Value = '  -1.45%  '
